# MAI_holdings.xlsx update
# - Bump the "as of" date in the confidential disclaimer text from 2021-05-17 to 2021-05-18
# - Refresh the Weight (D) / Percent Change (E) figures for rows 2-7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect (using the workbook's own password) so the
# cell values below can be updated, then restore protection afterwards.
$ws.Unprotect("D382")

# Update the disclaimer paragraph's "as of" date.
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-18 for illustrative purposes only and are subject to change."

# Row 2 (PFFD)
$ws.Range("D2").Value = 0.4750227568390289
$ws.Range("E2").Value = 0.0003930817610062753

# Row 3 (VYM)
$ws.Range("D3").Value = 0.3432849205384004
$ws.Range("E3").Value = -0.01004838109415707

# Row 4 (USRT)
$ws.Range("D4").Value = 0.09579370080868163
$ws.Range("E4").Value = 0.0008976660682225024

# Row 5 (HYG)
$ws.Range("D5").Value = 0.05346782287378645
$ws.Range("E5").Value = -0.002412406662837374

# Row 6 (AMLP)
$ws.Range("D6").Value = 0.03243079894010252
$ws.Range("E6").Value = -0.01505681818181803

# Row 7 (Total)
$ws.Range("E7").Value = -0.003794034944185953

# Restore sheet protection with the original password.
$ws.Protect("D382")
